$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.144.43"
$ws.Range("E2").Value = "  -3.43%  "

$ws.Range("D3").Value = "'1.849.99"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'0.7040"
$ws.Range("E5").Value = "  -5.44%  "

$ws.Range("D6").Value = "'238.19"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.3042"
$ws.Range("E8").Value = "  -4.21%  "

$ws.Range("D9").Value = "'0.07467"
$ws.Range("E9").Value = "  +3.18%  "

$ws.Range("E10").Value = "  -6.73%  "

$ws.Range("D11").Value = "'0.08124"

$ws.Range("D12").Value = "'0.7249"
$ws.Range("E12").Value = "  -5.02%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.829.51"
$ws.Range("E13").Value = "  -4.98%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.217"
$ws.Range("E14").Value = "  -4.23%  "

$ws.Range("D15").Value = "'88.99"
$ws.Range("E15").Value = "  -4.35%  "

$ws.Range("D16").Value = "'29.110.87"
$ws.Range("E16").Value = "  -3.75%  "

$ws.Range("D17").Value = "'5.777"
$ws.Range("E17").Value = "  -6.65%  "

$ws.Range("D18").Value = "'238.60"
$ws.Range("E18").Value = "  -4.88%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'13.06"
$ws.Range("E19").Value = "  -4.53%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007659"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'2.091.33"
$ws.Range("E22").Value = "  -4.74%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "'7.549"
$ws.Range("E24").Value = "  -5.89%  "

$ws.Range("D25").Value = "'162.02"
$ws.Range("E25").Value = "  -1.38%  "

$ws.Range("D26").Value = "'8.982"
$ws.Range("E26").Value = "  -3.57%  "

$ws.Range("D27").Value = "'0.1460"
$ws.Range("E27").Value = "  -8.00%  "

$ws.Range("D28").Value = "'17.99"
$ws.Range("E28").Value = "  -4.38%  "

$ws.Range("E29").Value = "  -6.55%  "

$ws.Range("D30").Value = "'1.386"
$ws.Range("E30").Value = "  -6.40%  "

$ws.Range("D31").Value = "'4.557"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").Value = "'1.493"
$ws.Range("E32").Value = "  -2.85%  "

$ws.Range("D33").Value = "'3.994"
$ws.Range("E33").Value = "  -5.56%  "

$ws.Range("D34").Value = "'0.05151"
$ws.Range("E34").Value = "  -4.74%  "

$ws.Range("D35").Value = "'1.187"
$ws.Range("E35").Value = "  -5.31%  "

$ws.Range("D36").Value = "'1.038"
$ws.Range("E36").Value = "  +4.40%  "

$ws.Range("D37").Value = "'0.7003"
$ws.Range("E37").Value = "  -8.87%  "

$ws.Range("E38").Value = "  -2.81%  "

$ws.Range("D39").Value = "'0.01867"
$ws.Range("E39").Value = "  -5.51%  "

$ws.Range("D40").Value = "'2.677"
$ws.Range("E40").Value = "  -3.40%  "

$ws.Range("D41").Value = "'0.9556"
$ws.Range("E41").Value = "  +9.72%  "

$ws.Range("D42").Value = "'6.003"
$ws.Range("E42").Value = "  -1.51%  "

$ws.Range("D43").Value = "'1.080.10"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D44").Value = "'0.4298"
$ws.Range("E44").Value = "  -6.05%  "

$ws.Range("D45").Value = "'69.84"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'102.21"
$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.749"
$ws.Range("E48").Value = "  -6.67%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'1.979.45"
$ws.Range("E49").Value = "  -4.08%  "

$ws.Range("D50").Value = "'9.164"
$ws.Range("E50").Value = "  -4.88%  "

$ws.Range("D51").Value = "'7.044"
$ws.Range("E51").Value = "  -7.74%  "

